$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 20 with forecast data for 2025 -> 2026
$ws.Range("A20").Value = 45986
$ws.Range("B20").Value = 2025
$ws.Range("C20").Value = 2.560577522109297
$ws.Range("D20").Value = 2026
$ws.Range("E20").Value = 2.991302072731838

# Copy the formatting (date style) from the row above (A19) into A20
# so the new date cell renders/serializes with the same style index.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
